# "recalculation - use of Ref numbers" test case
#
# The quantitative sheet's exchange table is replaced with a smaller,
# UUID-keyed data set that exercises Ref('exchange_id') / Ref('exchange_id',
# 'parameter_id') style mathematical relations. The meta sheet is untouched
# (its cell text is identical before/after - only shared-string bookkeeping
# changed in the source XML, which Excel handles on its own when we merely
# rewrite the quantitative sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quantitative")

$refFormula = "Ref('e22d5152-bfce-4032-8ef8-48a7be435495') + Ref('e22d5152-bfce-4032-8ef8-48a7be435495', 'ed3491b8-e2ad-47c8-be56-59e2f9b1deb1') + Ref('e22d5152-bfce-4032-8ef8-48a7be435495', 'ed3491b8-e2ad-47c8-be56-59e2f9b1deb1') + Ref('e22d5152-bfce-4032-8ef8-48a7be435495', Production Volume)"

# Drop the old rows 8-12 entirely (table shrinks from 12 rows to 7).
$ws.Range("A8:U12").ClearContents()

# --- Row 2: exchange (reference product electricity) ---
$ws.Range("A2").Value = "exchanges"
$ws.Range("B2").Value = "e22d5152-bfce-4032-8ef8-48a7be435495"
$ws.Range("E2").Value = "reference product"
$ws.Range("F2").Value = "electricity"
$ws.Range("J2").Value = "allocatable"
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = "kWh"
$ws.Range("N2").ClearContents()

# --- Row 3: property on that exchange ---
$ws.Range("A3").Value = "properties"
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = "4765611c-3ce3-4a1b-9ae7-0888225fb670"
$ws.Range("E3").Value = "reference product"
$ws.Range("F3").Value = "electricity"
$ws.Range("J3").Value = "allocatable"
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = "EURO2005"
$ws.Range("N3").ClearContents()

# --- Row 4: another property on that exchange ---
$ws.Range("A4").Value = "properties"
$ws.Range("B4").ClearContents()
$ws.Range("C4").Value = "ed3491b8-e2ad-47c8-be56-59e2f9b1deb1"
$ws.Range("E4").Value = "reference product"
$ws.Range("F4").Value = "electricity"
$ws.Range("J4").Value = "allocatable"
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = "kg"
$ws.Range("N4").ClearContents()

# --- Row 5: parameter "p" ---
$ws.Range("A5").Value = "parameter"
$ws.Range("B5").ClearContents()
$ws.Range("C5").Value = "ed3491b8-e2ad-47c8-be56-59e2f9b1deb2"
$ws.Range("E5").ClearContents()
$ws.Range("F5").Value = "p"
$ws.Range("G5").ClearContents()
$ws.Range("H5").ClearContents()
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = "kg"
$ws.Range("N5").ClearContents()

# --- Row 6: production volume exchange ---
$ws.Range("A6").Value = "production volume"
$ws.Range("B6").Value = "e22d5152-bfce-4032-8ef8-48a7be435495"
$ws.Range("E6").Value = "reference product"
$ws.Range("F6").Value = "electricity"
$ws.Range("G6").ClearContents()
$ws.Range("H6").ClearContents()
$ws.Range("J6").Value = "allocatable"
$ws.Range("K6").Value = 1000
$ws.Range("L6").Value = "kWh"
$ws.Range("N6").ClearContents()

# --- Row 7: parameter "q" = Ref math relation ---
$ws.Range("A7").Value = "parameter"
$ws.Range("B7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("F7").Value = "q"
$ws.Range("J7").ClearContents()
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = "kg"
$ws.Range("M7").Value = $refFormula
$ws.Range("N7").ClearContents()

# Match the saved selection/active cell (first empty row below the table).
$null = $ws.Range("K8").Select()
